$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that changed from 45171 (2023-09-02)
# to 45172 (2023-09-03) for every data row (rows 2 through 359).
$ws.Range("C2:C359").Value = 45172
